$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for handoff
#
# The "097298bd-7ea0-4fd0-9d70-83728cc19d14" file has dropped out of the
# localization report (it's handed-off set is gone), and the
# "05daeeb2-ac12-4594-84bc-3e5a63870673" file's status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", with refreshed
# handoff timestamps. This touches the Overview sheet plus the per-locale
# (zh-cn / de-de) detail sheets; each loses its 097298bd row and the row
# that used to follow it (.localization-config) shifts up.
# ---------------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# --- Overview sheet ---------------------------------------------------
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Rows.Item(3).Delete()

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/.localization-config", "", "", ".localization-config") | Out-Null

# --- zh-cn detail sheet -------------------------------------------------
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "2016-01-14 03:16:51"
$ws2.Rows.Item(3).Delete()

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8f2c89afa37d26fcbee6a2ba0b590fff2233bbd4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2c3f80dd30e368918e65bf793c13107ea8ad58cd/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ee874b779de1283ab0ad35ba5c148c0ae63a6934/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/.localization-config", "", "", ".localization-config") | Out-Null

# --- de-de detail sheet -------------------------------------------------
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "2016-01-14 03:17:02"
$ws3.Rows.Item(3).Delete()

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb4976e3f62c9ff00ed333f6dbc2d63dd696f792/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/05b017fe35d3376bad03bc64d3d2b7c1267f6469/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/101d8147b9135a4c1af57843dfb0fd15259f1f09/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b507373a6101f9fc256c88b1ce7ea876d0b91af5/.localization-config", "", "", ".localization-config") | Out-Null
